$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 4 through 8 (old iteration data no longer present)
$ws.Range("A4:C8").EntireRow.Delete()

# Update row 2 values (B2, C2)
$ws.Range("B2").Value = "[0;-1;1]"
$ws.Range("C2").Value = "inf"

# Update row 3 values (B3, C3)
$ws.Range("B3").Value = "[0;-1;1]"
$ws.Range("C3").Value = "nan"
